# Update workbook for "Add data for 2022-05-09"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-01"

# Update the header label in I1 (shared string "2022 (through 04-30)" -> "2022 (through 05-01)")
$ws.Range("I1").Value = "2022 (through 05-01)"

# Update changed numeric values
$ws.Range("I5").Value = 116
$ws.Range("I6").Value = 1
$ws.Range("H11").Value = 194
$ws.Range("H14").Value = 1851
